# Add a new "Verify defined asset ownership standard exists" requirement row
# to the "Vulnerability Management" worksheet, inserted above the existing
# row 9 ("2. Assets Types" / "Verify all business critical servers ...").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Vulnerability Management")

# Insert a new blank row above row 9; everything below shifts down by one.
$ws.Rows.Item(9).Insert()

# Populate the new row with the requirement text and Level 1/2/3 answers.
$ws.Range("B9").Value = "Verify defined asset ownership standard exists"
$ws.Range("C9").Value = "N"
$ws.Range("D9").Value = "Y"
$ws.Range("E9").Value = "Y"

# Match the style (centered alignment) used by the other answer cells.
$ws.Range("C9:E9").HorizontalAlignment = -4108

# Reflect the cursor position recorded in the saved workbook.
$ws.Range("E9").Select()
